# finished up flashcards for Module 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 23: the "Stress Management" question gets a more specific title.
# (Its answer in B23 stays the same.)
$ws.Range("A23").Value = "Stress Management of Pharmacy Technician"

# Append the new flashcards (rows 24-29).
$rows = @(
    @("Communication of Pharmacy Technician", "good verbal (attitude and tone) and non-verbal (open body language, eye contact) communication skills to pharmacists, patients, other staff, other healthcare professionals (doctors, nurses) (inter-professional)"),
    @("Cultural Competence of Pharmacy Technician", "Cultural and linguistic competence is a set of congruent behaviors, attitudes, and policies that come together in a system, agency, or among professionals that enable effective work in cross-cultural situations"),
    @("Cultural Competence of Pharmacy Technician (cont)", "able to effectively work with those that are different from them, ability to effectively treat those that are of a differnet race, gender, religion is essential"),
    @("Professionalism of Pharmacy Technician", "conduct, aims, or qualities that mark a professional person. Combination of appearance (good hygiene, well-groomed/dressed), personality (integrity, trustworhtiness), and actions (timeliness, active listening)"),
    @("How to get certified as a Pharmacy Technicain?", "Pharmacy Technician Certification exam (PTCE), or Exam for the Certification of Pharmacy Technicians (ExCPT)"),
    @("How to get a license as a Pharmacy Technician?", "Need to register with your state in order to practice")
)

$r = 24
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# Match the author's final selection state.
$ws.Range("B29").Select()
